# Update - Task Schedular.
# The automated test-data capture tool (Katalon) re-ran its supplier-credential
# workflow and stamped new timestamped values into row 2 of Sheet1. This
# script reproduces that same data refresh: the UserName/Password/
# NewPassword columns are untouched, while SupplierName, WorkGroupName,
# FormTask, DocTask, ACKTask, ItemName, FormTask1, DocTask1 and ACKTask1
# are updated with the latest run's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Supplier_05/03/19-15:02"
$ws.Range("D2").Value = "SuppReq_05/03/19-15:03"
$ws.Range("E2").Value = "WorkGroup_08/03/19-14:43"
$ws.Range("F2").Value = "TestFormTask_05/03/19-15:04"
$ws.Range("G2").Value = "TestDocTask_05/03/19-15:05"
$ws.Range("H2").Value = "TestACKTask_05/03/19-15:04"
$ws.Range("I2").Value = "ItemReq_05/03/19-15:06"
$ws.Range("J2").Value = "TestFormTask_05/03/19-15:07"
$ws.Range("K2").Value = "TestDocTask_05/03/19-15:08"
$ws.Range("L2").Value = "TestACKTask_05/03/19-15:06"

$wb.Save()
